$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from row 17 to row 18 first so values aren't reinterpreted
$ws.Range("A17:H17").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 18 data
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 43212
$ws.Range("C18").Value = 14.13
$ws.Range("D18").Value = 94.8
$ws.Range("E18").Value = 209
$ws.Range("F18").Formula = "=E18-E17"
$ws.Range("G18").Formula = "=ROUND((D18/1.88)/1.88,2)"
$ws.Range("H18").Value = 20.6

$ws.Range("F13").Select()
